# Fixing eeprom settings and naming conventions in controller
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the error code identifiers from camelCase to UPPER_SNAKE_CASE
$ws.Range("B2").Value = "ERR_INVALID_ADDRESS"
$ws.Range("B3").Value = "ERR_PACKET_LENGTH"
$ws.Range("B4").Value = "ERR_UNKOWN_COMMAND"
$ws.Range("B5").Value = "ERR_UNKOWN_PARAMETERS"
$ws.Range("B6").Value = "ERR_OUTSIDE_RANGE"
$ws.Range("B7").Value = "ERR_UNEXPECTED_PACKET"

# Add a new error code row for under-voltage detection
$ws.Range("A8").Value = 56
$ws.Range("B8").Value = "ERR_UNDERVOLTAGE"
$ws.Range("C8").Value = "The 5v supply voltage has dropped unexpectedly."
$ws.Range("D8").Value = "The supply voltage (in mV)"

# Move the active selection to the newly added cell
$ws.Range("D8").Select()
